$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(781, 1).Value = 9
$ws.Cells.Item(781, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(781, 3).Value = "Metropolitana"
$ws.Cells.Item(781, 4).Value = 44911
$ws.Cells.Item(781, 5).Value = 13
$ws.Cells.Item(781, 6).Value = "Fruta"
$ws.Cells.Item(781, 7).Value = 100103
$ws.Cells.Item(781, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(781, 9).Value = 100103006
$ws.Cells.Item(781, 10).Value = "Nectarín"
$ws.Cells.Item(781, 11).Value = "Artic Star"
$ws.Cells.Item(781, 12).Value = "Especial"
$ws.Cells.Item(781, 13).Value = 300
$ws.Cells.Item(781, 14).Value = 19200
$ws.Cells.Item(781, 15).Value = 19200
$ws.Cells.Item(781, 16).Value = 19200
$ws.Cells.Item(781, 17).Value = "`$/caja 16 kilos granel"
$ws.Cells.Item(781, 18).Value = "Región Metropolitana"
$ws.Cells.Item(781, 19).Value = 1200
$ws.Cells.Item(781, 20).Value = 16
$ws.Cells.Item(781, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(782, 1).Value = 9
$ws.Cells.Item(782, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(782, 3).Value = "Metropolitana"
$ws.Cells.Item(782, 4).Value = 44911
$ws.Cells.Item(782, 5).Value = 13
$ws.Cells.Item(782, 6).Value = "Fruta"
$ws.Cells.Item(782, 7).Value = 100103
$ws.Cells.Item(782, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(782, 9).Value = 100103006
$ws.Cells.Item(782, 10).Value = "Nectarín"
$ws.Cells.Item(782, 11).Value = "Artic Star"
$ws.Cells.Item(782, 12).Value = "Primera"
$ws.Cells.Item(782, 13).Value = 290
$ws.Cells.Item(782, 14).Value = 16000
$ws.Cells.Item(782, 15).Value = 16000
$ws.Cells.Item(782, 16).Value = 16000
$ws.Cells.Item(782, 17).Value = "`$/caja 16 kilos granel"
$ws.Cells.Item(782, 18).Value = "Región Metropolitana"
$ws.Cells.Item(782, 19).Value = 1000
$ws.Cells.Item(782, 20).Value = 16
$ws.Cells.Item(782, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(783, 1).Value = 9
$ws.Cells.Item(783, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(783, 3).Value = "Metropolitana"
$ws.Cells.Item(783, 4).Value = 44911
$ws.Cells.Item(783, 5).Value = 13
$ws.Cells.Item(783, 6).Value = "Fruta"
$ws.Cells.Item(783, 7).Value = 100103
$ws.Cells.Item(783, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(783, 9).Value = 100103006
$ws.Cells.Item(783, 10).Value = "Nectarín"
$ws.Cells.Item(783, 11).Value = "Artic Star"
$ws.Cells.Item(783, 12).Value = "Segunda"
$ws.Cells.Item(783, 13).Value = 380
$ws.Cells.Item(783, 14).Value = 12800
$ws.Cells.Item(783, 15).Value = 12800
$ws.Cells.Item(783, 16).Value = 12800
$ws.Cells.Item(783, 17).Value = "`$/caja 16 kilos granel"
$ws.Cells.Item(783, 18).Value = "Región Metropolitana"
$ws.Cells.Item(783, 19).Value = 800
$ws.Cells.Item(783, 20).Value = 16
$ws.Cells.Item(783, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(784, 1).Value = 9
$ws.Cells.Item(784, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(784, 3).Value = "Metropolitana"
$ws.Cells.Item(784, 4).Value = 44911
$ws.Cells.Item(784, 5).Value = 13
$ws.Cells.Item(784, 6).Value = "Fruta"
$ws.Cells.Item(784, 7).Value = 100103
$ws.Cells.Item(784, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(784, 9).Value = 100103006
$ws.Cells.Item(784, 10).Value = "Nectarín"
$ws.Cells.Item(784, 11).Value = "Early Glo"
$ws.Cells.Item(784, 12).Value = "Especial"
$ws.Cells.Item(784, 13).Value = 300
$ws.Cells.Item(784, 14).Value = 15000
$ws.Cells.Item(784, 15).Value = 15000
$ws.Cells.Item(784, 16).Value = 15000
$ws.Cells.Item(784, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(784, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(784, 19).Value = 1000
$ws.Cells.Item(784, 20).Value = 15
$ws.Cells.Item(784, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(785, 1).Value = 9
$ws.Cells.Item(785, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(785, 3).Value = "Metropolitana"
$ws.Cells.Item(785, 4).Value = 44911
$ws.Cells.Item(785, 5).Value = 13
$ws.Cells.Item(785, 6).Value = "Fruta"
$ws.Cells.Item(785, 7).Value = 100103
$ws.Cells.Item(785, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(785, 9).Value = 100103006
$ws.Cells.Item(785, 10).Value = "Nectarín"
$ws.Cells.Item(785, 11).Value = "Early Glo"
$ws.Cells.Item(785, 12).Value = "Primera"
$ws.Cells.Item(785, 13).Value = 350
$ws.Cells.Item(785, 14).Value = 12000
$ws.Cells.Item(785, 15).Value = 12000
$ws.Cells.Item(785, 16).Value = 12000
$ws.Cells.Item(785, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(785, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(785, 19).Value = 800
$ws.Cells.Item(785, 20).Value = 15
$ws.Cells.Item(785, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(786, 1).Value = 9
$ws.Cells.Item(786, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(786, 3).Value = "Metropolitana"
$ws.Cells.Item(786, 4).Value = 44911
$ws.Cells.Item(786, 5).Value = 13
$ws.Cells.Item(786, 6).Value = "Fruta"
$ws.Cells.Item(786, 7).Value = 100103
$ws.Cells.Item(786, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(786, 9).Value = 100103006
$ws.Cells.Item(786, 10).Value = "Nectarín"
$ws.Cells.Item(786, 11).Value = "Early John"
$ws.Cells.Item(786, 12).Value = "Especial"
$ws.Cells.Item(786, 13).Value = 290
$ws.Cells.Item(786, 14).Value = 18000
$ws.Cells.Item(786, 15).Value = 18000
$ws.Cells.Item(786, 16).Value = 18000
$ws.Cells.Item(786, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(786, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(786, 19).Value = 1000
$ws.Cells.Item(786, 20).Value = 18
$ws.Cells.Item(786, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(787, 1).Value = 9
$ws.Cells.Item(787, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(787, 3).Value = "Metropolitana"
$ws.Cells.Item(787, 4).Value = 44911
$ws.Cells.Item(787, 5).Value = 13
$ws.Cells.Item(787, 6).Value = "Fruta"
$ws.Cells.Item(787, 7).Value = 100103
$ws.Cells.Item(787, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(787, 9).Value = 100103006
$ws.Cells.Item(787, 10).Value = "Nectarín"
$ws.Cells.Item(787, 11).Value = "Early John"
$ws.Cells.Item(787, 12).Value = "Primera"
$ws.Cells.Item(787, 13).Value = 300
$ws.Cells.Item(787, 14).Value = 14400
$ws.Cells.Item(787, 15).Value = 14400
$ws.Cells.Item(787, 16).Value = 14400
$ws.Cells.Item(787, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(787, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(787, 19).Value = 800
$ws.Cells.Item(787, 20).Value = 18
$ws.Cells.Item(787, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(788, 1).Value = 9
$ws.Cells.Item(788, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(788, 3).Value = "Metropolitana"
$ws.Cells.Item(788, 4).Value = 44911
$ws.Cells.Item(788, 5).Value = 13
$ws.Cells.Item(788, 6).Value = "Fruta"
$ws.Cells.Item(788, 7).Value = 100103
$ws.Cells.Item(788, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(788, 9).Value = 100103006
$ws.Cells.Item(788, 10).Value = "Nectarín"
$ws.Cells.Item(788, 11).Value = "Early John"
$ws.Cells.Item(788, 12).Value = "Segunda"
$ws.Cells.Item(788, 13).Value = 280
$ws.Cells.Item(788, 14).Value = 10800
$ws.Cells.Item(788, 15).Value = 10800
$ws.Cells.Item(788, 16).Value = 10800
$ws.Cells.Item(788, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(788, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(788, 19).Value = 600
$ws.Cells.Item(788, 20).Value = 18
$ws.Cells.Item(788, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Added rows 781-788"
